$d = $word.ActiveDocument

# Remove the _GoBack bookmark up front so it does not interfere with
# range/insert operations; it will be re-created at the very end in its
# new (correct) location.
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

function Add-Chunk {
    param($pos, $text, $bold, $italic)
    $r = $d.Range($pos, $pos)
    $r.InsertAfter($text)
    $r.Font.Name = "Times New Roman"
    if ($bold) { $r.Font.Bold = 1 } else { $r.Font.Bold = 0 }
    if ($italic) { $r.Font.Italic = 1 } else { $r.Font.Italic = 0 }
    return $r.End
}

# --- Paragraph 0: append onto the existing "7. " paragraph ---
$p = $d.Paragraphs.Last
$pos = $p.Range.End - 1
$pos = Add-Chunk $pos "While the ECE students are correct that the " $false $false
$pos = Add-Chunk $pos "worst case" $true $false
$pos = Add-Chunk $pos " p" $false $false
$pos = Add-Chunk $pos "erformance of hash tables is O(" $false $false
$pos = Add-Chunk $pos "n" $false $true
$pos = Add-Chunk $pos "), they do not consider " $false $false
$pos = Add-Chunk $pos "the expected frequency of the worst case" $false $false
$pos = Add-Chunk $pos " nor the expected (average) performance of hash tables." $false $false
$pos = Add-Chunk $pos " For use in a real project, you must look at these factors as well as worst case performance to decide whether or not to use a certain data structure. Through analysis of these factors, I will show that the ECE students are wrong." $false $false

# --- Paragraph 1 ---
$p = $d.Paragraphs.Last
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$pos = $p.Range.Start

# --- Paragraph 2 ---
$p = $d.Paragraphs.Last
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$pos = $p.Range.Start
$pos = Add-Chunk $pos "To determine the average case performance of a linear probing hash table, you must address what happens in each possible case of operation. Assuming a dynamically resizing hash table, the two possible cases are that the hash table does or doesn’t need to be resized. If the table is resized, a new array must be allocated and every key must be re-hashed into the new array. This whole operation runs in O(" $false $false
$pos = Add-Chunk $pos "n" $false $true
$pos = Add-Chunk $pos ") time. However, this case occurs infrequently (only when λ ≥ ½)." $false $false

# --- Paragraph 3 ---
$p = $d.Paragraphs.Last
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$pos = $p.Range.Start

# --- Paragraph 4 ---
$p = $d.Paragraphs.Last
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$pos = $p.Range.Start
$pos = Add-Chunk $pos "Furthermore, this " $false $false
$pos = Add-Chunk $pos "“" $false $false
$pos = Add-Chunk $pos "heavy" $false $false
$pos = Add-Chunk $pos "”" $false $false
$pos = Add-Chunk $pos " cost can be distributed " $false $false
$pos = Add-Chunk $pos "across all the following operations before the next array resizing. When an operation occurs that doesn’t result in array resizing, it runs in constant time. In an amortized analysis where you distribute the present cost of resizing the hash table across future operations, it can be proven that the average performance " $false $false
$pos = Add-Chunk $pos "of the basic functions is O(1), constant time. " $false $false

# --- Paragraph 5 ---
$p = $d.Paragraphs.Last
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$pos = $p.Range.Start

# --- Paragraph 6 ---
$p = $d.Paragraphs.Last
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$pos = $p.Range.Start
$pos = Add-Chunk $pos "The flaw in the ECE students’ logic is that they don’t take into account the average performance. While the worst case is O(" $false $false
$pos = Add-Chunk $pos "n" $false $true
$pos = Add-Chunk $pos "), this case occurs infrequently and the expected performance of any arbitrary function call is O(1). Therefore, on average, the performance of these functions is O(1)." $false $false
$pos = Add-Chunk $pos " " $false $false
$pos = Add-Chunk $pos "Furthermore, you could simply initialize the hash table to be more than double the possible number of elements that will be inserted into it. Therefore, the table will never be resized and the performance of the specified functions will always be O(1). Thereby, you completely eliminate the possibility of the worst case occurring. " $false $false

# --- Paragraph 7 ---
$p = $d.Paragraphs.Last
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$pos = $p.Range.Start

# --- Paragraph 8 ---
$p = $d.Paragraphs.Last
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$pos = $p.Range.Start

# --- Paragraph 9 ---
$p = $d.Paragraphs.Last
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$pos = $p.Range.Start
$pos = Add-Chunk $pos "8." $false $false
$pos = Add-Chunk $pos " In a project where there will be any number of calls the RangeCount, I would suggest using a Balan" $false $false
$pos = Add-Chunk $pos "ced BST instead of a hash table. I would use a BBST because RangeCount’s worst case performance for a BBST is O(log " $false $false
$pos = Add-Chunk $pos "n" $false $true
$pos = Add-Chunk $pos ") where a Hash Table’s worst case is O(" $false $false
$pos = Add-Chunk $pos "h" $false $true
$pos = Add-Chunk $pos "). This difference is very significant." $false $false
$pos = Add-Chunk $pos " " $false $false
$pos = Add-Chunk $pos "As shown in Reports 1 & 2, the performance of RangeCount " $false $false
$pos = Add-Chunk $pos "in a BBST is orders of magnitude faster than Hash Tables, even when the BBST is operating on data sets orders of magnitude larger than the data sets in the Hash Table." $false $false
$pos = Add-Chunk $pos " For this reason, I would suggest that the CS 240 student use a BBST instead of a Hash Table." $false $false

# --- Paragraph 10 ---
$p = $d.Paragraphs.Last
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$pos = $p.Range.Start

# --- Paragraph 11 ---
$p = $d.Paragraphs.Last
$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$pos = $p.Range.Start
$pos = Add-Chunk $pos "9. " $false $false

# --- Re-create the _GoBack bookmark at the very end of the document ---
$endPos = $d.Content.End - 1
$bmRange = $d.Range($endPos, $endPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

Write-Output "done"